# Actualización desde MV -datos-
# Append new daily rate observations (10, 13, 14, 15, 16, 20 September 2021)
# to the end of the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Serie (col A), TPM (col B), Facilidad liquidez (col C), Facilidad deposito (col D)
$newRows = @(
    @("10-09-2021", 1.5, 1.75, 1.25),
    @("13-09-2021", 1.5, 1.75, 1.25),
    @("14-09-2021", 1.5, 1.75, 1.25),
    @("15-09-2021", 1.5, 1.75, 1.25),
    @("16-09-2021", 1.5, 1.75, 1.25),
    @("20-09-2021", 1.5, $null, $null)
)

# Locate the first empty row right after the current data block (row 175 -> 176)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$r = $lastRow + 1

foreach ($row in $newRows) {
    $cellA = $ws.Cells.Item($r, 1)
    # Force the "dd-mm-yyyy" looking label to stay plain text (as the rest of
    # column A already is) instead of being auto-parsed into a date serial,
    # the same way typing a leading apostrophe in Excel would; then restore
    # the default cell style so the cell itself carries no extra formatting.
    $cellA.Value = "'" + $row[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($null -ne $row[2]) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    if ($null -ne $row[3]) {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
    $r = $r + 1
}
